$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D2").Value = "55.067.89"

# Row 3
$ws.Range("E3").Value = "  -5.28%  "
$ws.Range("D3").Value = "2.341.07"

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("E5").Value = "  -2.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "475.46"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.26"
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("E7").Value = "  +22.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("D7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "  -5.50%  "
$ws.Range("D9").Value = "2.341.88"

# Row 10
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0963"
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "  -6.31%  "

# Row 12
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.323"
$ws.Range("D12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "  +1.13%  "

# Row 14
$ws.Range("E14").Value = "  -5.45%  "
$ws.Range("D14").Value = "2.747.29"

# Row 15
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D15").Value = "55.055.54"

# Row 16
$ws.Range("E16").Value = "  -5.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.96"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "  -4.38%  "

# Row 18
$ws.Range("E18").Value = "  -5.55%  "
$ws.Range("D18").Value = "2.340.49"

# Row 19
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.56"
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "314.45"
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "  -4.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.59"
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "  -2.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.61"
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.74"
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "  -4.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.393"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "  -6.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.152"
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "  -5.29%  "
$ws.Range("D28").Value = "2.443.87"

# Row 29
$ws.Range("E29").Value = "  -7.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.02"
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "  +0.15%  "

# Row 31
$ws.Range("E31").Value = "  -5.44%  "
$ws.Range("D31").Value = "0.0₃0745"

# Row 32
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.16"
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "  -3.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "144.37"
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "  -2.27%  "

# Row 35
$ws.Range("E35").Value = "  -1.59%  "

# Row 36
$ws.Range("E36").Value = "  -4.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.57"
$ws.Range("D36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "  -4.68%  "

# Row 38
$ws.Range("E38").Value = "  -5.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.812"
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E39").Value = "  +8.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.65"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.40"
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.32"
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = "  -4.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.574"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "  -6.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0518"
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.16"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "249.94"
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0219"
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "  -8.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.33"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "  -5.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.65"
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "  -5.35%  "
$ws.Range("D51").Value = "1.769.23"

